$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (R / 1K 0603 Resistor): quantity on board 2 -> 3, note updated to "buying 4"
$ws.Range("D5").Value = 3
$ws.Range("I5").Value = "Buying 10 is cheaper than buying 4."

# Row 9 (1uF Ceramic Cap): quantity on board 3 -> 4, order 5 -> 6
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 6

# Update the active selection to C13
$ws.Range("C13").Select()

$wb.Save()
